$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.348.01"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").Value = "1.846.71"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9972"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6271"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9986"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07607"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2901"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.31%  "
$ws.Range("E10").Value = "  +0.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07728"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.024"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6786"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.00001047"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.93"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.141"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("D17").Value = "29.368.64"
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "227.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.33"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9982"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9981"
$ws.Range("D22").Style = "Normal"
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("E24").Value = "  -0.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.406"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.52%  "
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.401"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.457"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05600"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.106"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.064"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.35%  "
$ws.Range("E32").Value = "  +0.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.834"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.6965"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.586"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("D37").Value = "1.227.55"
$ws.Range("E37").Value = "  -0.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.718"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.352"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9008"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9987"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "101.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.63"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.208"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.53%  "
$ws.Range("E45").Value = "  -1.58%  "
$ws.Range("E46").Value = "  -0.57%  "
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("E48").Value = "  -0.75%  "
$ws.Range("E49").Value = "  +1.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05699"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4619"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.17%  "
